{"js": "// This script updates the \"date of practice\" heading and the 25\n// three-digit-by-one-digit division problems in the table below it.\n// Because every old value is unique within the document, each one can\n// be located with a simple body search and swapped for its replacement\n// without touching any other run's formatting.\nconst replacements = [\n  [\"2025-12-08 Monday\", \"2025-12-09 Tuesday\"],\n  [\"262\u00f75=52, 2\", \"971\u00f78=121, 3\"],\n  [\"460\u00f72=230, 0\", \"119\u00f75=23, 4\"],\n  [\"766\u00f75=153, 1\", \"824\u00f72=412, 0\"],\n  [\"750\u00f75=150, 0\", \"283\u00f76=47, 1\"],\n  [\"793\u00f74=198, 1\", \"173\u00f73=57, 2\"],\n  [\"857\u00f76=142, 5\", \"205\u00f74=51, 1\"],\n  [\"486\u00f79=54, 0\", \"607\u00f78=75, 7\"],\n  [\"561\u00f76=93, 3\", \"571\u00f79=63, 4\"],\n  [\"803\u00f78=100, 3\", \"176\u00f73=58, 2\"],\n  [\"414\u00f74=103, 2\", \"377\u00f72=188, 1\"],\n  [\"239\u00f76=39, 5\", \"551\u00f75=110, 1\"],\n  [\"293\u00f72=146, 1\", \"332\u00f72=166, 0\"],\n  [\"547\u00f75=109, 2\", \"947\u00f75=189, 2\"],\n  [\"904\u00f79=100, 4\", \"678\u00f72=339, 0\"],\n  [\"117\u00f77=16, 5\", \"392\u00f73=130, 2\"],\n  [\"291\u00f74=72, 3\", \"860\u00f75=172, 0\"],\n  [\"974\u00f77=139, 1\", \"587\u00f79=65, 2\"],\n  [\"199\u00f77=28, 3\", \"224\u00f78=28, 0\"],\n  [\"450\u00f74=112, 2\", \"245\u00f77=35, 0\"],\n  [\"835\u00f76=139, 1\", \"777\u00f74=194, 1\"],\n  [\"470\u00f74=117, 2\", \"538\u00f77=76, 6\"],\n  [\"356\u00f79=39, 5\", \"679\u00f73=226, 1\"],\n  [\"279\u00f73=93, 0\", \"724\u00f79=80, 4\"],\n  [\"615\u00f79=68, 3\", \"193\u00f76=32, 1\"],\n  [\"884\u00f72=442, 0\", \"388\u00f72=194, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Updates the \"date of practice\" heading and the 25 three-digit-by-\n# one-digit division problems in the table below it. Every \"old\" value\n# below is unique within the document, so Find/Replace can locate and\n# swap each one without disturbing any other run's formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-12-08 Monday\", \"2025-12-09 Tuesday\"),\n  @(\"262\u00f75=52, 2\", \"971\u00f78=121, 3\"),\n  @(\"460\u00f72=230, 0\", \"119\u00f75=23, 4\"),\n  @(\"766\u00f75=153, 1\", \"824\u00f72=412, 0\"),\n  @(\"750\u00f75=150, 0\", \"283\u00f76=47, 1\"),\n  @(\"793\u00f74=198, 1\", \"173\u00f73=57, 2\"),\n  @(\"857\u00f76=142, 5\", \"205\u00f74=51, 1\"),\n  @(\"486\u00f79=54, 0\", \"607\u00f78=75, 7\"),\n  @(\"561\u00f76=93, 3\", \"571\u00f79=63, 4\"),\n  @(\"803\u00f78=100, 3\", \"176\u00f73=58, 2\"),\n  @(\"414\u00f74=103, 2\", \"377\u00f72=188, 1\"),\n  @(\"239\u00f76=39, 5\", \"551\u00f75=110, 1\"),\n  @(\"293\u00f72=146, 1\", \"332\u00f72=166, 0\"),\n  @(\"547\u00f75=109, 2\", \"947\u00f75=189, 2\"),\n  @(\"904\u00f79=100, 4\", \"678\u00f72=339, 0\"),\n  @(\"117\u00f77=16, 5\", \"392\u00f73=130, 2\"),\n  @(\"291\u00f74=72, 3\", \"860\u00f75=172, 0\"),\n  @(\"974\u00f77=139, 1\", \"587\u00f79=65, 2\"),\n  @(\"199\u00f77=28, 3\", \"224\u00f78=28, 0\"),\n  @(\"450\u00f74=112, 2\", \"245\u00f77=35, 0\"),\n  @(\"835\u00f76=139, 1\", \"777\u00f74=194, 1\"),\n  @(\"470\u00f74=117, 2\", \"538\u00f77=76, 6\"),\n  @(\"356\u00f79=39, 5\", \"679\u00f73=226, 1\"),\n  @(\"279\u00f73=93, 0\", \"724\u00f79=80, 4\"),\n  @(\"615\u00f79=68, 3\", \"193\u00f76=32, 1\"),\n  @(\"884\u00f72=442, 0\", \"388\u00f72=194, 0\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  # MatchCase=$true, Forward=$true, Wrap=wdFindContinue(1), Replace=wdReplaceAll(2)\n  $find.Execute([ref]$oldText, $true, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2) | Out-Null\n}\n"}
